# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
# Values that look numeric are entered with a leading apostrophe so Excel
# stores them as text (matching the source data's text-typed Price column)
# instead of silently coercing them to a float.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.258.08"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "2.567.72"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'584.46"
$ws.Range("E5").Value = "  +2.86%  "

$ws.Range("D6").Value = "'149.07"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +0.89%  "

$ws.Range("E9").Value = "  +3.58%  "

$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "'27.88"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").Value = "3.026.96"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").Value = "63.137.69"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "'0.0000147"
$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("D17").Value = "2.569.53"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("D18").Value = "'11.46"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").Value = "'340.70"
$ws.Range("E19").Value = "  +1.48%  "

$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("D21").Value = "'6.87"
$ws.Range("E21").Value = "  +1.60%  "

$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "'66.31"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").Value = "2.689.81"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("D27").Value = "'8.22"
$ws.Range("E27").Value = "  +14.53%  "

$ws.Range("D28").Value = "'8.56"
$ws.Range("E28").Value = "  +1.96%  "

$ws.Range("E29").Value = "  +1.33%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").Value = "'1.97"
$ws.Range("E31").Value = "  +5.67%  "

$ws.Range("E32").Value = "  +2.16%  "

$ws.Range("D33").Value = "'177.70"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").Value = "'441.63"
$ws.Range("E34").Value = "  +6.90%  "

$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("E36").Value = "  +2.37%  "

$ws.Range("D37").Value = "'19.37"
$ws.Range("E37").Value = "  +2.75%  "

$ws.Range("D38").Value = "'4.52"
$ws.Range("E38").Value = "  +2.97%  "

$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").Value = "'152.07"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").Value = "'3.84"
$ws.Range("E43").Value = "  +2.78%  "

$ws.Range("D44").Value = "'21.50"
$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("D45").Value = "'0.0553"
$ws.Range("E45").Value = "  +5.90%  "

$ws.Range("D46").Value = "'0.609"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("E47").Value = "  +1.35%  "

$ws.Range("E48").Value = "  +3.35%  "

$ws.Range("D49").Value = "'18.53"
$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("D50").Value = "'1.75"
$ws.Range("E50").Value = "  -1.44%  "

$ws.Range("E51").Value = "  -0.24%  "

Write-Output "Applied 79 cell updates to cryptos sheet"
